$wb = $excel.ActiveWorkbook

# --- Rename the "Type" worksheet to "CreatureType" ---
$wsType = $wb.Worksheets.Item("Type")
$wsType.Name = "CreatureType"

# --- Creature sheet: rename column header "Size" -> "SizeType" ---
$wsCreature = $wb.Worksheets.Item("Creature")
$wsCreature.Range("F1").Value = "SizeType"

# --- CreatureInstance sheet: rename column header "Creature (FK)" -> "CreatureID (FK)" ---
$wsCreatureInstance = $wb.Worksheets.Item("CreatureInstance")
$wsCreatureInstance.Range("E1").Value = "CreatureID (FK)"

# --- Update selections/active cells to match the edited session ---
$wsCreature.Range("C20").Select()
$wsCreatureInstance.Range("E1").Select()
$wsType.Range("C32").Select()

# Leave TypeImmunity as the final active sheet/tab (matches activeTab retained in workbook view)
$wsTypeImmunity = $wb.Worksheets.Item("TypeImmunity")
$wsTypeImmunity.Select()
